$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数 / interested count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 97
$ws1.Range("F3").Value = 4067
$ws1.Range("F4").Value = 2374
$ws1.Range("F5").Value = 472
$ws1.Range("F6").Value = 16
$ws1.Range("F8").Value = 30
$ws1.Range("F9").Value = 199
$ws1.Range("F11").Value = 84
$ws1.Range("F12").Value = 136
$ws1.Range("F13").Value = 1522
$ws1.Range("F14").Value = 274
$ws1.Range("F15").Value = 2931
$ws1.Range("F16").Value = 201

# Sheet "全部类型" (All types) - update column F (想去人数 / interested count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 97
$ws4.Range("F3").Value = 4067
$ws4.Range("F4").Value = 2374
$ws4.Range("F5").Value = 472
$ws4.Range("F6").Value = 16
$ws4.Range("F9").Value = 30
$ws4.Range("F11").Value = 199
$ws4.Range("F13").Value = 84
$ws4.Range("F14").Value = 136
$ws4.Range("F17").Value = 1522
$ws4.Range("F18").Value = 274
$ws4.Range("F19").Value = 2931
$ws4.Range("F20").Value = 201
